$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (386-464), continuing the series through 2021-12-08
$newRows = @(
    @(386, 44460, 0, 5, 57.49770009199631),
    @(387, 44461, 1, 6, 68.99724011039559),
    @(388, 44462, 0, 5, 57.49770009199631),
    @(389, 44463, 2, 6, 68.99724011039559),
    @(390, 44464, 1, 6, 68.99724011039559),
    @(391, 44465, 3, 8, 91.99632014719411),
    @(392, 44466, 5, 12, 137.9944802207912),
    @(393, 44467, 1, 13, 149.4940202391904),
    @(394, 44468, 0, 12, 137.9944802207912),
    @(395, 44469, 9, 21, 241.4903403863846),
    @(396, 44470, 1, 20, 229.9908003679853),
    @(397, 44471, 4, 23, 264.489420423183),
    @(398, 44472, 0, 20, 229.9908003679853),
    @(399, 44473, 2, 17, 195.4921803127875),
    @(400, 44474, 0, 16, 183.9926402943882),
    @(401, 44475, 0, 16, 183.9926402943882),
    @(402, 44476, 2, 9, 103.4958601655934),
    @(403, 44477, 1, 9, 103.4958601655934),
    @(404, 44478, 0, 5, 57.49770009199631),
    @(405, 44479, 0, 5, 57.49770009199631),
    @(406, 44480, 1, 4, 45.99816007359706),
    @(407, 44481, 0, 4, 45.99816007359706),
    @(408, 44482, 0, 4, 45.99816007359706),
    @(409, 44483, 0, 2, 22.99908003679853),
    @(410, 44484, 0, 1, 11.49954001839926),
    @(411, 44485, 0, 1, 11.49954001839926),
    @(412, 44486, 0, 1, 11.49954001839926),
    @(413, 44487, 0, 0, 0),
    @(414, 44488, 0, 0, 0),
    @(415, 44489, 0, 0, 0),
    @(416, 44490, 0, 0, 0),
    @(417, 44491, 0, 0, 0),
    @(418, 44492, 0, 0, 0),
    @(419, 44493, 0, 0, 0),
    @(420, 44494, 1, 1, 11.49954001839926),
    @(421, 44495, 0, 1, 11.49954001839926),
    @(422, 44496, 0, 1, 11.49954001839926),
    @(423, 44497, 1, 2, 22.99908003679853),
    @(424, 44498, 0, 2, 22.99908003679853),
    @(425, 44499, 0, 2, 22.99908003679853),
    @(426, 44500, 1, 3, 34.4986200551978),
    @(427, 44501, 0, 2, 22.99908003679853),
    @(428, 44502, 0, 2, 22.99908003679853),
    @(429, 44503, 0, 2, 22.99908003679853),
    @(430, 44504, 0, 1, 11.49954001839926),
    @(431, 44505, 0, 1, 11.49954001839926),
    @(432, 44506, 2, 3, 34.4986200551978),
    @(433, 44507, 0, 2, 22.99908003679853),
    @(434, 44508, 0, 2, 22.99908003679853),
    @(435, 44509, 0, 2, 22.99908003679853),
    @(436, 44510, 0, 2, 22.99908003679853),
    @(437, 44511, 2, 4, 45.99816007359706),
    @(438, 44512, 0, 4, 45.99816007359706),
    @(439, 44513, 0, 2, 22.99908003679853),
    @(440, 44514, 0, 2, 22.99908003679853),
    @(441, 44515, 0, 2, 22.99908003679853),
    @(442, 44516, 12, 14, 160.9935602575897),
    @(443, 44517, 0, 14, 160.9935602575897),
    @(444, 44518, 1, 13, 149.4940202391904),
    @(445, 44519, 0, 13, 149.4940202391904),
    @(446, 44520, 1, 14, 160.9935602575897),
    @(447, 44521, 0, 14, 160.9935602575897),
    @(448, 44522, 2, 16, 183.9926402943882),
    @(449, 44523, 2, 6, 68.99724011039559),
    @(450, 44524, 6, 12, 137.9944802207912),
    @(451, 44525, 0, 11, 126.4949402023919),
    @(452, 44526, 8, 19, 218.491260349586),
    @(453, 44527, 1, 19, 218.491260349586),
    @(454, 44528, 0, 19, 218.491260349586),
    @(455, 44529, 1, 18, 206.9917203311867),
    @(456, 44530, 1, 17, 195.4921803127875),
    @(457, 44531, 0, 11, 126.4949402023919),
    @(458, 44532, 2, 13, 149.4940202391904),
    @(459, 44533, 3, 8, 91.99632014719411),
    @(460, 44534, 4, 11, 126.4949402023919),
    @(461, 44535, 1, 12, 137.9944802207912),
    @(462, 44536, 5, 16, 183.9926402943882),
    @(463, 44537, 2, 17, 195.4921803127875),
    @(464, 44538, 1, 18, 206.9917203311867)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Copy the style of the preceding row's date cell (column A) so the new
    # date cell keeps the same number format / border style (xf index 2)
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
